$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D2: "Pilt" label, styled like the other header cells (bold Arial 10,
#         wrap text) but with a left+right "medium grey" border instead of
#         the full box border used by the other header cells. ---
$d2 = $ws.Range("D2")

# Start from the existing header style (A2/B2/C2) so font/border/fill/
# alignment line up, then only touch what needs to change.
$ws.Range("A2").Copy($d2)
$d2.Value2 = "Pilt"

# Drop the top & bottom edges of the copied border so only left/right remain.
$d2.Borders.Item(8).LineStyle = -4142
$d2.Borders.Item(9).LineStyle = -4142

# --- D6: plain "example.jpg" text, default styling (same as D1's note). ---
$ws.Range("D6").Value2 = "example.jpg"

# --- Selection moves to D7 (matches the saved cursor position). ---
$ws.Range("D7").Select()

Write-Output "done"
